# Penalty Reward System (unfinished) — remove stale/incorrect weekly and
# monthly data points from the PO analysis sheets.
#
# "Weekly Quantity": the rows for 2023-06-11, 2023-06-18 and 2023-07-16
# (old rows 5-7) are removed; later rows shift up so the sheet ends at B7.
#
# "Monthly Trend": the rows for 2023-07 and 2023-08-early (old rows 3-4)
# are removed; the remaining row shifts up so the sheet ends at B3.

$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Rows("5:7").Delete()

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Rows("3:4").Delete()
